$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a literal text value into a cell without Excel's
# autodetection turning ambiguous day/month strings (e.g. "01-08-2022")
# into a date serial number, and without leaving behind any NumberFormat /
# quote-prefix style residue. We do this by writing a text-returning
# formula (="...") and then collapsing it down to its static value via
# Copy + PasteSpecial(xlPasteValues = -4163).
function Set-TextValue($addr, $text) {
    $r = $ws.Range($addr)
    $r.Formula = "=""" + $text + """"
    $r.Copy()
    $r.PasteSpecial(-4163)
}

# Column A: date strings, slashes -> hyphens
Set-TextValue "A3"  "28-07-2022"
Set-TextValue "A4"  "01-08-2022"
Set-TextValue "A5"  "04-08-2022"
Set-TextValue "A6"  "08-08-2022"
Set-TextValue "A7"  "11-08-2022"
Set-TextValue "A8"  "15-08-2022"
Set-TextValue "A9"  "18-08-2022"
Set-TextValue "A10" "22-08-2022"
Set-TextValue "A11" "25-08-2022"
Set-TextValue "A12" "29-08-2022"
Set-TextValue "A13" "01-09-2022"
Set-TextValue "A14" "05-09-2022"
Set-TextValue "A15" "08-09-2022"
Set-TextValue "A16" "12-09-2022"
Set-TextValue "A17" "15-09-2022"
Set-TextValue "A18" "19-09-2022"
Set-TextValue "A19" "22-09-2022"
Set-TextValue "A20" "26-09-2022"
Set-TextValue "A21" "29-09-2022"

# Updated attendance counts
$ws.Range("D3").Value  = 1
$ws.Range("G3").Value  = 1

$ws.Range("D10").Value = 1
$ws.Range("G10").Value = 1

$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0

$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 1
$ws.Range("H21").Value = 0
